$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Heading3 title
Replace-Text "Environmental Impacts and adaptation" "Environmental Impacts Assessment"

# Metadata list (Créditos-trabalho / Carga horária / Ativação).
# These three lines live in three separate, identically-formatted <w:r> runs
# inside the SAME paragraph. Touching that paragraph with more than one
# mutating call causes the engine to coalesce adjacent same-format runs, so
# do all three substitutions in a single Find/Execute call (one run split on
# manual line breaks, matched literally via the vertical-tab run-break char).
$brk = [char]11
$metaOld = "Créditos-trabalho: 0" + $brk + "Carga horária: 60 h" + $brk + "Ativação: 01/01/2018"
$metaNew = "Créditos-trabalho: 1" + $brk + "Carga horária: 90 h" + $brk + "Ativação: 01/01/2022"
Replace-Text $metaOld $metaNew

# Objetivos (PT)
Replace-Text "Propiciar os conhecimentos básicos sobre o uso dos recursos ambientais pelas atividades antrópicas, as interferências no meio, os impactos associados e o encaminhamento para adequação ambiental, alinhado aos princípios da sustentabilidade. Conhecer ferramentas de apoio estratégico na identificação de problemas e estabelecimento de ações de correção." "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."

# Objetivos (EN, italic)
Replace-Text "Provide basic knowledge on the use of environmental resources by anthropic activities, environmental interference, associated impacts and referral to environmental compliance, in line with the principles of sustainability. To know support tools in identifying problems and establishing corrective actions." "Provide knowledge on the fundamentals, objectives and methods of environmental impact assessment."

# Programa resumido (PT)
Replace-Text "Processo de avaliação de impacto ambiental e seus objetivos; Impactos ambientais sobre águas superficiais, subterrâneas e oceânicas; Poluição hídrica; Impactos ambientais sobre o solo; Poluição atmosférica e saúde humana." "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."

# Programa (PT)
Replace-Text "Disponibilidade hídrica; Impactos promovidos por reservatórios e sua classificação; Fontes de poluição hídrica; Fontes de contaminação de águas subterrâneas; Ocupação desordenada da faixa litorânea; Degradação dos solos; Impactos de poluentes atmosféricos que afetam a saúde humana; Métodos científicos utilizados nos estados de poluição; Impactos das mudanças climáticas; Identificação de impactos; previsão de impactos; avaliação da importância dos impactos; Etapa de planejamento e da elaboração de um estudo de impacto ambiental; Estudo de Impacto Ambiental (EIA); Relatório de Impacto Ambiental (RIMA)." "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."

# Programa (EN, italic)
Replace-Text "Water availability; Impacts promoted by reservoirs and their classification; Sources of water pollution; Sources of groundwater contamination; Disordered occupation of the coastal strip; Soil degradation; Impacts of air pollutants that affect human health; Scientific methods used in pollution states; Impacts of climate change; Identification of impacts; Prediction of impacts; Assessment of the importance of impacts; Stage of planning and preparation of an environmental impact study; Environmental Impact Study ; Environmental Impact Report." "Basic concepts and definitions. Origin and dissemination of the Environmental Impact Assessment (EIA). EIA and licensing: objectives and fundamentals. Brazilian legal and institutional framework for EIA. The EIA process and its components. Steps in planning and execution of an Environmental Impact Statement. Technological and localization alternatives. Baseline studies. Impact identification and prediction techniques. Methods and criteria for determining impact significance. Environmental management plan: mitigation, compensation, monitoring. Decision making and follow-up. Case studies."

# Avaliação: Método / Critério / Norma de recuperação
Replace-Text "Aula expositiva, exercícios dirigidos, viagens para aula prática, trabalho prático e seminários." "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
Replace-Text "Média ponderada de exercícios e provas." "Média ponderada de atividades e provas."
Replace-Text "Prova única com nota igual ou superior a 5,0." "1 (uma) prova escrita"

# Bibliografia: collapse the multi-run / multi-<w:br/> paragraph into a single run of text.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Bibliografia b")) {
        $targetPara = $p
        break
    }
}
if ($targetPara -ne $null) {
    $r = $d.Range($targetPara.Range.Start, $targetPara.Range.End - 1)
    $r.Text = "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
}
